# Case 2_82 (380 kV slack setpoint change 1.05 -> 1.02 p.u.) results refresh.
# res_bus/vm_pu sheet: rewrite the recalculated per-unit bus voltage magnitudes
# (columns B-F and I-N, rows 2-25) with the values produced by the new power flow run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02  # B2
$ws.Cells.Item(2, 3).Value = 1.040005016769034  # C2
$ws.Cells.Item(2, 4).Value = 1.047274939994667  # D2
$ws.Cells.Item(2, 5).Value = 1.05417882408134  # E2
$ws.Cells.Item(2, 6).Value = 1.060185003763752  # F2
$ws.Cells.Item(2, 9).Value = 1.036047759160601  # I2
$ws.Cells.Item(2, 10).Value = 1.045093964709403  # J2
$ws.Cells.Item(2, 11).Value = 1.050038033861895  # K2
$ws.Cells.Item(2, 12).Value = 1.056922749147953  # L2
$ws.Cells.Item(2, 13).Value = 1.062912472444702  # M2
$ws.Cells.Item(2, 14).Value = 1.018857188523417  # N2
$ws.Cells.Item(3, 2).Value = 1.02  # B3
$ws.Cells.Item(3, 3).Value = 1.041227385148773  # C3
$ws.Cells.Item(3, 4).Value = 1.048229230044684  # D3
$ws.Cells.Item(3, 5).Value = 1.055287085791918  # E3
$ws.Cells.Item(3, 6).Value = 1.06132604761513  # F3
$ws.Cells.Item(3, 9).Value = 1.036274064803214  # I3
$ws.Cells.Item(3, 10).Value = 1.045960409227279  # J3
$ws.Cells.Item(3, 11).Value = 1.050803765620319  # K3
$ws.Cells.Item(3, 12).Value = 1.057843448092422  # L3
$ws.Cells.Item(3, 13).Value = 1.063867067013442  # M3
$ws.Cells.Item(3, 14).Value = 1.01915179761539  # N3
$ws.Cells.Item(4, 2).Value = 1.02  # B4
$ws.Cells.Item(4, 3).Value = 1.042018060085472  # C4
$ws.Cells.Item(4, 4).Value = 1.04884631720711  # D4
$ws.Cells.Item(4, 5).Value = 1.056004763385076  # E4
$ws.Cells.Item(4, 6).Value = 1.062064627941135  # F4
$ws.Cells.Item(4, 9).Value = 1.036418952961264  # I4
$ws.Cells.Item(4, 10).Value = 1.046520276846107  # J4
$ws.Cells.Item(4, 11).Value = 1.051298223163536  # K4
$ws.Cells.Item(4, 12).Value = 1.058439168959357  # L4
$ws.Cells.Item(4, 13).Value = 1.064484416970963  # M4
$ws.Cells.Item(4, 14).Value = 1.019341987831093  # N4
$ws.Cells.Item(5, 2).Value = 1.02  # B5
$ws.Cells.Item(5, 3).Value = 1.04235039498187  # C5
$ws.Cells.Item(5, 4).Value = 1.049105644986671  # D5
$ws.Cells.Item(5, 5).Value = 1.056306609365805  # E5
$ws.Cells.Item(5, 6).Value = 1.062375187523558  # F5
$ws.Cells.Item(5, 9).Value = 1.036479493817172  # I5
$ws.Cells.Item(5, 10).Value = 1.046755459497151  # J5
$ws.Cells.Item(5, 11).Value = 1.051505849112641  # K5
$ws.Cells.Item(5, 12).Value = 1.058689602789207  # L5
$ws.Cells.Item(5, 13).Value = 1.064743871110539  # M5
$ws.Cells.Item(5, 14).Value = 1.019421838323568  # N5
$ws.Cells.Item(6, 2).Value = 1.02  # B6
$ws.Cells.Item(6, 3).Value = 1.042406191728215  # C6
$ws.Cells.Item(6, 4).Value = 1.049149181671052  # D6
$ws.Cells.Item(6, 5).Value = 1.056357298531003  # E6
$ws.Cells.Item(6, 6).Value = 1.062427335412344  # F6
$ws.Cells.Item(6, 9).Value = 1.036489637193651  # I6
$ws.Cells.Item(6, 10).Value = 1.046794936839395  # J6
$ws.Cells.Item(6, 11).Value = 1.051540696134775  # K6
$ws.Cells.Item(6, 12).Value = 1.058731651326755  # L6
$ws.Cells.Item(6, 13).Value = 1.064787429946692  # M6
$ws.Cells.Item(6, 14).Value = 1.019435239394439  # N6
$ws.Cells.Item(7, 2).Value = 1.02  # B7
$ws.Cells.Item(7, 3).Value = 1.042022501012953  # C7
$ws.Cells.Item(7, 4).Value = 1.048849782732683  # D7
$ws.Cells.Item(7, 5).Value = 1.056008796137187  # E7
$ws.Cells.Item(7, 6).Value = 1.062068777414991  # F7
$ws.Cells.Item(7, 9).Value = 1.036419763365284  # I7
$ws.Cells.Item(7, 10).Value = 1.046523420096594  # J7
$ws.Cells.Item(7, 11).Value = 1.051300998429538  # K7
$ws.Cells.Item(7, 12).Value = 1.058442515296011  # L7
$ws.Cells.Item(7, 13).Value = 1.064487884122385  # M7
$ws.Cells.Item(7, 14).Value = 1.019343055211429  # N7
$ws.Cells.Item(8, 2).Value = 1.02  # B8
$ws.Cells.Item(8, 3).Value = 1.040418180876752  # C8
$ws.Cells.Item(8, 4).Value = 1.04759753047943  # D8
$ws.Cells.Item(8, 5).Value = 1.054553251164569  # E8
$ws.Cells.Item(8, 6).Value = 1.060570573923641  # F8
$ws.Cells.Item(8, 9).Value = 1.036124560570736  # I8
$ws.Cells.Item(8, 10).Value = 1.045386945793136  # J8
$ws.Cells.Item(8, 11).Value = 1.050297028546953  # K8
$ws.Cells.Item(8, 12).Value = 1.057233910878808  # L8
$ws.Cells.Item(8, 13).Value = 1.063235152312523  # M8
$ws.Cells.Item(8, 14).Value = 1.01895684461846  # N8
$ws.Cells.Item(9, 2).Value = 1.02  # B9
$ws.Cells.Item(9, 3).Value = 1.037588936976311  # C9
$ws.Cells.Item(9, 4).Value = 1.04538779100996  # D9
$ws.Cells.Item(9, 5).Value = 1.051992638097964  # E9
$ws.Cells.Item(9, 6).Value = 1.057932412188114  # F9
$ws.Cells.Item(9, 9).Value = 1.03559252184424  # I9
$ws.Cells.Item(9, 10).Value = 1.043378313526491  # J9
$ws.Cells.Item(9, 11).Value = 1.048520045680992  # K9
$ws.Cells.Item(9, 12).Value = 1.055103913474594  # L9
$ws.Cells.Item(9, 13).Value = 1.061025065241877  # M9
$ws.Cells.Item(9, 14).Value = 1.01827289944554  # N9
$ws.Cells.Item(10, 2).Value = 1.02  # B10
$ws.Cells.Item(10, 3).Value = 1.035701150477749  # C10
$ws.Cells.Item(10, 4).Value = 1.043912495107848  # D10
$ws.Cells.Item(10, 5).Value = 1.050288368435459  # E10
$ws.Cells.Item(10, 6).Value = 1.056174832016505  # F10
$ws.Cells.Item(10, 9).Value = 1.03522984763189  # I10
$ws.Cells.Item(10, 10).Value = 1.042035111481236  # J10
$ws.Cells.Item(10, 11).Value = 1.047330059449694  # K10
$ws.Cells.Item(10, 12).Value = 1.053683681268777  # L10
$ws.Cells.Item(10, 13).Value = 1.059549863759841  # M10
$ws.Cells.Item(10, 14).Value = 1.017814638618178  # N10
$ws.Cells.Item(11, 2).Value = 1.02  # B11
$ws.Cells.Item(11, 3).Value = 1.034883303986935  # C11
$ws.Cells.Item(11, 4).Value = 1.04327315815474  # D11
$ws.Cells.Item(11, 5).Value = 1.049551054963406  # E11
$ws.Cells.Item(11, 6).Value = 1.055414049609636  # F11
$ws.Cells.Item(11, 9).Value = 1.0350709094366  # I11
$ws.Cells.Item(11, 10).Value = 1.041452498793939  # J11
$ws.Cells.Item(11, 11).Value = 1.046813506675209  # K11
$ws.Cells.Item(11, 12).Value = 1.053068638903608  # L11
$ws.Cells.Item(11, 13).Value = 1.058910643291003  # M11
$ws.Cells.Item(11, 14).Value = 1.01761565775997  # N11
$ws.Cells.Item(12, 2).Value = 1.02  # B12
$ws.Cells.Item(12, 3).Value = 1.034579453322213  # C12
$ws.Cells.Item(12, 4).Value = 1.043035600160282  # D12
$ws.Cells.Item(12, 5).Value = 1.049277279968602  # E12
$ws.Cells.Item(12, 6).Value = 1.055131499105788  # F12
$ws.Cells.Item(12, 9).Value = 1.035011587248509  # I12
$ws.Cells.Item(12, 10).Value = 1.041235939218685  # J12
$ws.Cells.Item(12, 11).Value = 1.046621442473159  # K12
$ws.Cells.Item(12, 12).Value = 1.052840172889112  # L12
$ws.Cells.Item(12, 13).Value = 1.058673139778157  # M12
$ws.Cells.Item(12, 14).Value = 1.017541664260941  # N12
$ws.Cells.Item(13, 2).Value = 1.02  # B13
$ws.Cells.Item(13, 3).Value = 1.034644633352414  # C13
$ws.Cells.Item(13, 4).Value = 1.043086560783723  # D13
$ws.Cells.Item(13, 5).Value = 1.04933600130488  # E13
$ws.Cells.Item(13, 6).Value = 1.055192105451631  # F13
$ws.Cells.Item(13, 9).Value = 1.035024324982507  # I13
$ws.Cells.Item(13, 10).Value = 1.041282398851789  # J13
$ws.Cells.Item(13, 11).Value = 1.046662649669991  # K13
$ws.Cells.Item(13, 12).Value = 1.052889180179332  # L13
$ws.Cells.Item(13, 13).Value = 1.058724088209936  # M13
$ws.Cells.Item(13, 14).Value = 1.017557539886486  # N13
$ws.Cells.Item(14, 2).Value = 1.02  # B14
$ws.Cells.Item(14, 3).Value = 1.03485818895888  # C14
$ws.Cells.Item(14, 4).Value = 1.043253523170438  # D14
$ws.Cells.Item(14, 5).Value = 1.049528422683632  # E14
$ws.Cells.Item(14, 6).Value = 1.055390693130195  # F14
$ws.Cells.Item(14, 9).Value = 1.035066011671368  # I14
$ws.Cells.Item(14, 10).Value = 1.041434601006374  # J14
$ws.Cells.Item(14, 11).Value = 1.046797634541908  # K14
$ws.Cells.Item(14, 12).Value = 1.053049754080102  # L14
$ws.Cells.Item(14, 13).Value = 1.058891012582228  # M14
$ws.Cells.Item(14, 14).Value = 1.017609543132688  # N14
$ws.Cells.Item(15, 2).Value = 1.02  # B15
$ws.Cells.Item(15, 3).Value = 1.034989758784611  # C15
$ws.Cells.Item(15, 4).Value = 1.043356383603887  # D15
$ws.Cells.Item(15, 5).Value = 1.049646992542048  # E15
$ws.Cells.Item(15, 6).Value = 1.05551305454321  # F15
$ws.Cells.Item(15, 9).Value = 1.03509165837927  # I15
$ws.Cells.Item(15, 10).Value = 1.041528357687467  # J15
$ws.Cells.Item(15, 11).Value = 1.046880777503745  # K15
$ws.Cells.Item(15, 12).Value = 1.053148687353424  # L15
$ws.Cells.Item(15, 13).Value = 1.058993851075449  # M15
$ws.Cells.Item(15, 14).Value = 1.017641573019499  # N15
$ws.Cells.Item(16, 2).Value = 1.02  # B16
$ws.Cells.Item(16, 3).Value = 1.035755419133685  # C16
$ws.Cells.Item(16, 4).Value = 1.043954914664783  # D16
$ws.Cells.Item(16, 5).Value = 1.050337315032269  # E16
$ws.Cells.Item(16, 6).Value = 1.05622532803447  # F16
$ws.Cells.Item(16, 9).Value = 1.035240355822229  # I16
$ws.Cells.Item(16, 10).Value = 1.042073756430641  # J16
$ws.Cells.Item(16, 11).Value = 1.047364314278683  # K16
$ws.Cells.Item(16, 12).Value = 1.05372449799238  # L16
$ws.Cells.Item(16, 13).Value = 1.059592277166584  # M16
$ws.Cells.Item(16, 14).Value = 1.017827832683052  # N16
$ws.Cells.Item(17, 2).Value = 1.02  # B17
$ws.Cells.Item(17, 3).Value = 1.036235582704702  # C17
$ws.Cells.Item(17, 4).Value = 1.044330216335932  # D17
$ws.Cells.Item(17, 5).Value = 1.050770508429087  # E17
$ws.Cells.Item(17, 6).Value = 1.056672187434263  # F17
$ws.Cells.Item(17, 9).Value = 1.035333121452157  # I17
$ws.Cells.Item(17, 10).Value = 1.042415602619305  # J17
$ws.Cells.Item(17, 11).Value = 1.047667280443942  # K17
$ws.Cells.Item(17, 12).Value = 1.054085668684248  # L17
$ws.Cells.Item(17, 13).Value = 1.059967533131132  # M17
$ws.Cells.Item(17, 14).Value = 1.017944520688985  # N17
$ws.Cells.Item(18, 2).Value = 1.02  # B18
$ws.Cells.Item(18, 3).Value = 1.036515613405288  # C18
$ws.Cells.Item(18, 4).Value = 1.044549072955959  # D18
$ws.Cells.Item(18, 5).Value = 1.051023245264248  # E18
$ws.Cells.Item(18, 6).Value = 1.056932858274328  # F18
$ws.Cells.Item(18, 9).Value = 1.03538704698092  # I18
$ws.Cells.Item(18, 10).Value = 1.042614899587756  # J18
$ws.Cells.Item(18, 11).Value = 1.047843871968823  # K18
$ws.Cells.Item(18, 12).Value = 1.054296326517711  # L18
$ws.Cells.Item(18, 13).Value = 1.060186370395124  # M18
$ws.Cells.Item(18, 14).Value = 1.018012529701581  # N18
$ws.Cells.Item(19, 2).Value = 1.02  # B19
$ws.Cells.Item(19, 3).Value = 1.036611089831146  # C19
$ws.Cells.Item(19, 4).Value = 1.044623688873801  # D19
$ws.Cells.Item(19, 5).Value = 1.051109432619359  # E19
$ws.Cells.Item(19, 6).Value = 1.057021744627136  # F19
$ws.Cells.Item(19, 9).Value = 1.035405403146635  # I19
$ws.Cells.Item(19, 10).Value = 1.042682838441595  # J19
$ws.Cells.Item(19, 11).Value = 1.047904064196932  # K19
$ws.Cells.Item(19, 12).Value = 1.054368154251269  # L19
$ws.Cells.Item(19, 13).Value = 1.060260980956313  # M19
$ws.Cells.Item(19, 14).Value = 1.018035710017474  # N19
$ws.Cells.Item(20, 2).Value = 1.02  # B20
$ws.Cells.Item(20, 3).Value = 1.03618406991578  # C20
$ws.Cells.Item(20, 4).Value = 1.044289955249838  # D20
$ws.Cells.Item(20, 5).Value = 1.050724024425938  # E20
$ws.Cells.Item(20, 6).Value = 1.056624241042056  # F20
$ws.Cells.Item(20, 9).Value = 1.035323187518982  # I20
$ws.Cells.Item(20, 10).Value = 1.042378935712949  # J20
$ws.Cells.Item(20, 11).Value = 1.047634787825732  # K20
$ws.Cells.Item(20, 12).Value = 1.054046919211067  # L20
$ws.Cells.Item(20, 13).Value = 1.059927276194501  # M20
$ws.Cells.Item(20, 14).Value = 1.017932006668528  # N20
$ws.Cells.Item(21, 2).Value = 1.02  # B21
$ws.Cells.Item(21, 3).Value = 1.034795304000019  # C21
$ws.Cells.Item(21, 4).Value = 1.043204359119819  # D21
$ws.Cells.Item(21, 5).Value = 1.049471756740582  # E21
$ws.Cells.Item(21, 6).Value = 1.055332212971364  # F21
$ws.Cells.Item(21, 9).Value = 1.035053743863439  # I21
$ws.Cells.Item(21, 10).Value = 1.041389785453142  # J21
$ws.Cells.Item(21, 11).Value = 1.046757890211713  # K21
$ws.Cells.Item(21, 12).Value = 1.053002469405031  # L21
$ws.Cells.Item(21, 13).Value = 1.058841859417491  # M21
$ws.Cells.Item(21, 14).Value = 1.017594231771226  # N21
$ws.Cells.Item(22, 2).Value = 1.02  # B22
$ws.Cells.Item(22, 3).Value = 1.033921748239296  # C22
$ws.Cells.Item(22, 4).Value = 1.042521339620467  # D22
$ws.Cells.Item(22, 5).Value = 1.04868496156738  # E22
$ws.Cells.Item(22, 6).Value = 1.054520082587913  # F22
$ws.Cells.Item(22, 9).Value = 1.034882682161008  # I22
$ws.Cells.Item(22, 10).Value = 1.040766991433863  # J22
$ws.Cells.Item(22, 11).Value = 1.046205429803899  # K22
$ws.Cells.Item(22, 12).Value = 1.052345713031846  # L22
$ws.Cells.Item(22, 13).Value = 1.058159017496151  # M22
$ws.Cells.Item(22, 14).Value = 1.017381378023127  # N22
$ws.Cells.Item(23, 2).Value = 1.02  # B23
$ws.Cells.Item(23, 3).Value = 1.034384873895472  # C23
$ws.Cells.Item(23, 4).Value = 1.042883465335081  # D23
$ws.Cells.Item(23, 5).Value = 1.049102004307754  # E23
$ws.Cells.Item(23, 6).Value = 1.054950587958212  # F23
$ws.Cells.Item(23, 9).Value = 1.034973521887444  # I23
$ws.Cells.Item(23, 10).Value = 1.041097229844665  # J23
$ws.Cells.Item(23, 11).Value = 1.046498406045805  # K23
$ws.Cells.Item(23, 12).Value = 1.052693878921584  # L23
$ws.Cells.Item(23, 13).Value = 1.058521042977524  # M23
$ws.Cells.Item(23, 14).Value = 1.017494261569129  # N23
$ws.Cells.Item(24, 2).Value = 1.02  # B24
$ws.Cells.Item(24, 3).Value = 1.036207346453166  # C24
$ws.Cells.Item(24, 4).Value = 1.044308147657772  # D24
$ws.Cells.Item(24, 5).Value = 1.05074502835325  # E24
$ws.Cells.Item(24, 6).Value = 1.056645905874786  # F24
$ws.Cells.Item(24, 9).Value = 1.035327676801561  # I24
$ws.Cells.Item(24, 10).Value = 1.04239550420671  # J24
$ws.Cells.Item(24, 11).Value = 1.047649470222995  # K24
$ws.Cells.Item(24, 12).Value = 1.054064428451547  # L24
$ws.Cells.Item(24, 13).Value = 1.059945466704887  # M24
$ws.Cells.Item(24, 14).Value = 1.017937661379906  # N24
$ws.Cells.Item(25, 2).Value = 1.02  # B25
$ws.Cells.Item(25, 3).Value = 1.038320641261822  # C25
$ws.Cells.Item(25, 4).Value = 1.045959434867842  # D25
$ws.Cells.Item(25, 5).Value = 1.052654119956867  # E25
$ws.Cells.Item(25, 6).Value = 1.058614225377607  # F25
$ws.Cells.Item(25, 9).Value = 1.035731472556645  # I25
$ws.Cells.Item(25, 10).Value = 1.043898313297091  # J25
$ws.Cells.Item(25, 11).Value = 1.048980374534731  # K25
$ws.Cells.Item(25, 12).Value = 1.055654606828534  # L25
$ws.Cells.Item(25, 13).Value = 1.061596740874494  # M25
$ws.Cells.Item(25, 14).Value = 1.018450119393197  # N25
